$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the DD_object (C) and SD_object (E) values for rows 5-8
for ($r = 5; $r -le 8; $r++) {
    $cVal = $ws.Cells.Item($r, 3).Value2
    $eVal = $ws.Cells.Item($r, 5).Value2
    $ws.Cells.Item($r, 3).Value2 = $eVal
    $ws.Cells.Item($r, 5).Value2 = $cVal
}

# Update the active selection to match the latest working area
$ws.Range("E5:E8").Select() | Out-Null
